$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row r => @(A, C, D, E, F, G, H)   (B is always "falling")
$rows = @{
  2 = @(0, -1.181560516357422, 1.240549564361572, -0.5883253216743469, 0.008163382383910248, -0.07067988338795567, -0.05231227912008755)
  3 = @(100, -1.229650020599365, 1.595798969268799, -0.7779194116592407, -0.06663984263485132, -0.06482113017277281, -0.02257424863902029)
  4 = @(200, -0.2790107727050781, 1.227867126464844, -1.016065120697022, -0.005289537662809526, -0.05883742868900297, 0.1704730540513992)
  5 = @(300, -1.14784574508667, 1.330484867095947, -0.8695551156997681, 0.00837163017554712, -0.05375613881783047, 0.2080968008799986)
  6 = @(400, -1.288122177124023, 1.133394718170166, -0.7701815962791443, -0.05617183616215521, 0.01428591663187196, 0.0969609718092462)
  7 = @(500, -1.425555229187011, 1.231297016143799, -0.6422767043113708, 0.0280998013913631, 0.0762054398655891, -0.008552113547921099)
  8 = @(600, -1.529905319213867, 1.149116039276123, -0.9077204465866088, 0.06225272437388239, 0.09189357202161434, -0.06575131636451573)
  9 = @(700, -1.152251720428467, 1.056320667266846, -0.8245069980621338, 0.04992435195229265, 0.03760986301031981, -0.06972193650223989)
  10 = @(800, -1.123115062713623, 1.085736274719239, -0.7493376731872559, 0.03176499361341648, -0.02544809200546947, 0.0004164989699015492)
  11 = @(900, -1.394256114959717, 1.270269870758057, -0.8571128845214844, 0.03686016493222929, -0.01995030405338506, 0.01384165091440073)
  12 = @(1000, -1.345842838287354, 1.367433071136475, -0.7304041385650635, 0.01492454890500421, 0.02186619833281094, 0.002901610000249519)
  13 = @(1100, -1.356554985046387, 1.367568492889404, -0.6971900463104248, -0.02151911604133509, 0.04574547572569409, -0.03826237769151855)
  14 = @(1200, -1.285661697387695, 1.301629066467285, -0.5075737237930298, -0.01338350129398431, 0.03186217543076381, -0.08599316531961611)
  15 = @(1300, -1.321344375610352, 1.155209541320801, -0.5763433575630188, 0.02753058617765249, 0.01243943788788525, -0.08503521572459828)
  16 = @(1400, -1.26480770111084, 1.153616428375244, -0.6878960132598877, 0.02497605843977489, 0.07183220271359811, -0.08203642476688736)
  17 = @(1500, -0.8195595741271973, 1.118548393249511, -0.9026113748550416, 0.02700301970947864, 0.2729595926674929, -0.1162587641315027)
  18 = @(1600, -1.078746795654297, 1.204861640930176, -0.6440812349319458, 0.0572686158120632, 0.3712533414363861, -0.058643065392971)
  19 = @(1700, -1.571155071258545, 1.288459777832031, -0.5561246871948242, -0.1334879455918612, 0.2170098776167092, 0.09671107679605455)
  20 = @(1800, -1.56222677230835, 1.121483325958252, -0.1114475727081298, -0.05985091050917481, 0.1304891434582799, -0.121909264813769)
  21 = @(1900, -1.141444206237793, 1.122160911560059, 0.2474624365568161, -0.05091006410392845, 0.1894376386295665, -0.5621209686452691)
  22 = @(2000, -1.447634696960449, 0.8090605735778809, 2.089949369430542, 0.0685279762203039, -0.7451306256380917, -1.788169232281768)
  23 = @(2100, -2.48740816116333, -0.6178178787231445, 2.327399969100952, -0.745783193544906, -2.762346484444355, -1.756598580967298)
  24 = @(2200, -2.30517578125, -0.4891290664672851, 2.537726402282715, -0.1823294054378186, -0.2232989397915999, 1.138486255298958)
  25 = @(2300, -3.007162570953369, -1.502803325653076, 3.441044807434082, 0.5607048923319088, 2.322676355188556, 2.414278323000138)
  26 = @(2400, -3.378014087677002, -3.04712963104248, 5.235954284667969, -0.7038278525525871, -0.9540881785479458, -0.1982396515932949)
  27 = @(2500, -6.576094627380371, 0.4697372913360595, 2.406523942947388, -0.2033903666517955, -0.4635634205558086, 0.1531466895883727)
  28 = @(2600, 8.52888298034668, 7.144853115081787, -10.8126277923584, 0.07737163657491758, -0.3377668315714036, 0.5267046527429067)
  29 = @(2700, -10.79609775543213, -1.548766374588013, -0.7228314876556396, 0.437226802110672, -1.076497316360474, 0.305127203464508)
  30 = @(2800, -0.0029077529907226, 0.2878659963607788, 0.002394676208496, 0.3028364533727822, -0.5765596194700771, -0.2633939818902439)
  31 = @(2900, -3.11443567276001, 0.3510211706161499, -0.4145855903625488, -0.01123159040104305, -0.1864527369087388, -0.3065016323869878)
}

foreach ($r in ($rows.Keys | Sort-Object)) {
  $v = $rows[$r]
  $ws.Cells.Item($r, 1).Value = $v[0]
  $ws.Cells.Item($r, 2).Value = "falling"
  $ws.Cells.Item($r, 3).Value = $v[1]
  $ws.Cells.Item($r, 4).Value = $v[2]
  $ws.Cells.Item($r, 5).Value = $v[3]
  $ws.Cells.Item($r, 6).Value = $v[4]
  $ws.Cells.Item($r, 7).Value = $v[5]
  $ws.Cells.Item($r, 8).Value = $v[6]
}

Write-Output "Updated rows 2-31 with new falling sensor data"
